$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V2 Command")
$ws.Rows.Item(50).Insert()
$src = $ws.Range("B51:K51")
$dst = $ws.Range("B50:K50")
$src.Copy()
$dst.PasteSpecial(-4122)
